$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update designator list for the 0-ohm 1206 resistor row (row 6): add R28, R29
$ws.Range("B6").Value = "R4,R17,R15,R10,R21,R5,R3,R16,R24,R26,R28,R29"

# Update designator list for the 3k6 resistor row (row 12): add R9, R27
$ws.Range("B12").Value = "R2,R1,R9,R27"

# Row 13 changed from a 120R resistor to a 1k resistor, with an updated JLCPCB part number
$ws.Range("A13").Value = "1k"
$ws.Range("D13").Value = "C4410"

# Widen the Designator column to fit the longer designator lists
$ws.Columns.Item(2).ColumnWidth = 72.83

# Move the active selection to D13 (where edits were made)
$ws.Range("D13").Select()
